$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "test_existing_survey_import 1"

# Update row 2 - question type/code/text/detail
$ws.Range("B2").Value = "FreeText"
$ws.Range("A2").Value = "fdfuu42a22321c123a8_test"
$ws.Range("C2").Value = "Test question fdfuu42a22321c123a8_test"
$ws.Range("D2").Value = "Test question fdfuu42a22321c123a8_test"

# Update row 3 - question type/code/text/detail
$ws.Range("B3").Value = "FreeText"
$ws.Range("A3").Value = "fdfzz42a66321c123a8_test"
$ws.Range("C3").Value = "Test question fdfzz42a66321c123a8_test"
$ws.Range("D3").Value = "Test question fdfzz42a66321c123a8_test"

# Clear row 4 entirely (was a third test question, now removed)
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""

# Update the active selection to match the saved view state
$ws.Range("F15").Select()
